$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 06:46"

# Update per-country case statistics (daily refresh of COVID-19 data)
$ws.Range("B6").Value = 26878
$ws.Range("C6").Value = 2671
$ws.Range("E6").Value = 26352
$ws.Range("E23").Value = 1060
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 9
$ws.Range("B33").Value = 599
$ws.Range("C33").Value = 188
$ws.Range("D33").Value = 44
$ws.Range("E33").Value = 554
$ws.Range("B34").Value = 537
$ws.Range("D34").Value = 8
$ws.Range("E34").Value = 528
$ws.Range("F34").Value = 7
$ws.Range("H34").Value = 1
$ws.Range("B35").Value = 536
$ws.Range("D35").Value = 13
$ws.Range("E35").Value = 518
$ws.Range("F35").Value = 3
$ws.Range("H35").Value = 5
$ws.Range("B36").Value = 532
$ws.Range("D36").Value = 3
$ws.Range("E36").Value = 522
$ws.Range("F36").Value = 2
$ws.Range("H36").Value = 7
$ws.Range("B37").Value = 530
$ws.Range("D37").Value = 19
$ws.Range("E37").Value = 498
$ws.Range("F37").Value = 18
$ws.Range("H37").Value = 13
$ws.Range("B38").Value = 523
$ws.Range("D38").Value = 10
$ws.Range("E38").Value = 512
$ws.Range("F38").Value = 2
$ws.Range("H38").Value = 1
$ws.Range("B39").Value = 481
$ws.Range("D39").Value = 27
$ws.Range("E39").Value = 454
$ws.Range("F39").Value = 6
$ws.Range("H39").Value = 0
$ws.Range("B40").Value = 473
$ws.Range("D40").Value = 5
$ws.Range("E40").Value = 467
$ws.Range("F40").Value = 1
$ws.Range("H40").Value = 1
$ws.Range("B41").Value = 450
$ws.Range("D41").Value = 20
$ws.Range("E41").Value = 392
$ws.Range("F41").Value = 0
$ws.Range("H41").Value = 38
$ws.Range("B42").Value = 432
$ws.Range("D42").Value = 140
$ws.Range("E42").Value = 290
$ws.Range("F42").Value = 14
$ws.Range("H42").Value = 2
$ws.Range("B75").Value = 105
$ws.Range("C75").Value = 6
$ws.Range("D75").Value = 1
$ws.Range("E75").Value = 103
$ws.Range("F75").Value = 1
$ws.Range("H75").Value = 1
$ws.Range("B76").Value = 103
$ws.Range("D76").Value = 7
$ws.Range("E76").Value = 92
$ws.Range("F76").Value = 6
$ws.Range("H76").Value = 4
$ws.Range("B77").Value = 100
$ws.Range("D77").Value = 1
$ws.Range("E77").Value = 99
$ws.Range("F77").Value = 0
$ws.Range("H77").Value = 0

Write-Host "Update complete"
